$wb = $excel.ActiveWorkbook

# --- Sheet "Enums": insert a row above the enum table, turn it into a real
#     Excel Table named "Difficulty", and rename its header from "Type"
#     to "Difficulty" ------------------------------------------------------
$wsEnums = $wb.Worksheets.Item("Enums")

# Shift the Type/Value enum block from B2:C5 down to B3:C6
[void]$wsEnums.Rows("2:2").Insert()

# The old block carried hand-rolled header/zebra-stripe formatting; the
# edited sheet uses the default (unformatted) style instead.
[void]$wsEnums.Range("B3:C6").ClearFormats()

# Slightly wider first column
$wsEnums.Columns("B:B").ColumnWidth = 9.4

# Rename the header cell and promote the range to a real ListObject/Table
$wsEnums.Range("B3").Value = "Difficulty"
$loDifficulty = $wsEnums.ListObjects.Add(1, $wsEnums.Range("B3:C6"), 0, 1)
$loDifficulty.Name = "Difficulty"
$loDifficulty.TableStyle = "TableStyleMedium2"

# --- Sheet "Tables": re-point a few Difficulty cells -----------------------
$wsTables = $wb.Worksheets.Item("Tables")
$wsTables.Range("E4").Value = "Easy"
$wsTables.Range("E7").Value = "Medium"
$wsTables.Range("E10").Value = "Hard"

[void]$wsTables.Range("E10").Select()

# Leave "Enums" as the active / selected sheet, matching the saved file
[void]$wsEnums.Activate()
[void]$wsEnums.Range("D4").Select()
